# Update "Horarios Línea 141" workbook with the new scrape (07:12:53).
# Sheet "LP1912": header timestamps/row-count updated, rows 45-56 shuffled
#   into their new sorted-by-arrival-time positions, and 7 brand-new rows
#   (57-63) appended.
# Sheet "LP1912-215": header timestamps/row-count updated, 2 new rows
#   (15-16) appended.
# Sheet "6203-6173": only the "Última actualización" timestamp changes.

$wb = $excel.ActiveWorkbook

$newScrapTime = "07:12:53"

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newScrapTime"
$ws1.Range("A3").Value = "Total filas: 58"

$sheet1Rows = @(
    @(45, "07:12:53", "07:59", "23_HERNANDEZ", 47),
    @(46, "06:16:15", "08:00", "11_ETCHEVERRY", 104),
    @(47, "06:33:46", "08:00", "23_HERNANDEZ", 87),
    @(48, "06:45:50", "08:01", "23_HERNANDEZ", 76),
    @(49, "07:12:53", "08:01", "16_SANTA ANA", 49),
    @(50, "06:52:52", "08:02", "23_HERNANDEZ", 70),
    @(51, "06:16:15", "08:03", "17_ROMERO", 107),
    @(52, "06:16:15", "08:12", "10_OLMOS", 99),
    @(53, "07:12:53", "08:13", "10_OLMOS", 61),
    @(54, "06:16:15", "08:15", "17_ROMERO", 119),
    @(55, "06:33:46", "08:26", "15X38_ABASTO", 113),
    @(56, "06:33:46", "08:27", "84_COLONIA URQUIZA-ESC 49", 114),
    @(57, "06:45:50", "08:29", "14_ABASTO", 104),
    @(58, "06:33:46", "08:31", "16_P MOR-SANTA ANA", 118),
    @(59, "06:45:50", "08:38", "215C_EL PATO", 113),
    @(60, "07:12:53", "08:43", "10_OLMOS", 91),
    @(61, "07:12:53", "08:49", "215A_EL PATO", 97),
    @(62, "07:12:53", "08:59", "215B_EL PATO", 107),
    @(63, "07:12:53", "09:02", "17X38_ROMERO", 110)
)

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = "LP1912"
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newScrapTime"
$ws2.Range("A3").Value = "Total filas: 11"

$sheet2NewRows = @(
    @(15, "07:12:53", "08:49", "215A_EL PATO", 97),
    @(16, "07:12:53", "08:59", "215B_EL PATO", 107)
)

foreach ($row in $sheet2NewRows) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 5).Value = "LP1912"
}

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173 (only the header timestamp changes)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newScrapTime"
